# Insert two new offer rows ("Yuno Energy - D Smart Bonus + 6%" and
# "Energia Offer") above "Bord Gáis - Smart EV Bonus", pushing all
# subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4 and below down by two rows to make room for the new data.
$ws.Rows("4:5").Insert()

# Populate the newly inserted rows.
$ws.Range("A4").Value = "Yuno Energy - D Smart Bonus + 6%"
$ws.Range("B4").Value = 1402.35

$ws.Range("A5").Value = "Energia Offer"
$ws.Range("B5").Value = 1441.77
